$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: add column D "ppfd" ---
$ws.Range("D1").Value = "ppfd"

# --- Update existing rows 2-16: new A/B values, add D=800 ---
$ws.Range("A2").Value = 0.17062589280990501
$ws.Range("B2").Value = 13.1993016242526
$ws.Range("D2").Value = 800
$ws.Range("A3").Value = 5.0777736627691601
$ws.Range("B3").Value = 17.120787259933302
$ws.Range("D3").Value = 800
$ws.Range("A4").Value = 10.2560711073488
$ws.Range("B4").Value = 23.03793450082
$ws.Range("D4").Value = 800
$ws.Range("A5").Value = 15.290196285910699
$ws.Range("B5").Value = 25.357388497963001
$ws.Range("D5").Value = 800
$ws.Range("A6").Value = 20.173535791757001
$ws.Range("B6").Value = 22.0792550658695
$ws.Range("D6").Value = 800
$ws.Range("A7").Value = 0.14681762869689499
$ws.Range("B7").Value = 5.9996825564784704
$ws.Range("D7").Value = 800
$ws.Range("A8").Value = 5.1769747632400396
$ws.Range("B8").Value = 7.1192000423257804
$ws.Range("D8").Value = 800
$ws.Range("A9").Value = 10.2031638537643
$ws.Range("B9").Value = 7.0387810168773903
$ws.Range("D9").Value = 800
$ws.Range("A10").Value = 15.183720438071999
$ws.Range("B10").Value = -6.84090788847152
$ws.Range("D10").Value = 800
$ws.Range("A11").Value = 20.206602825247298
$ws.Range("B11").Value = -7.9212740066663301
$ws.Range("D11").Value = 800
$ws.Range("A12").Value = 0.0370350775091274
$ws.Range("B12").Value = -7.1985609227024998
$ws.Range("D12").Value = 800
$ws.Range("A13").Value = 5.1194381249669298
$ws.Range("B13").Value = -10.279879371461799
$ws.Range("D13").Value = 800
$ws.Range("A14").Value = 10.193243743717201
$ws.Range("B14").Value = -15.9610602613618
$ws.Range("D14").Value = 800
$ws.Range("A15").Value = 15.1658642399873
$ws.Range("B15").Value = -32.240622189302101
$ws.Range("D15").Value = 800
$ws.Range("A16").Value = 20.133194010898801
$ws.Range("B16").Value = -30.120099465636699
$ws.Range("D16").Value = 800

# --- Add new rows 17-31 (ppfd = 1600 block) ---
$ws.Range("A17").Value = -0.033500837520939797
$ws.Range("A17").NumberFormat = "0"
$ws.Range("B17").Value = 15.8793969849246
$ws.Range("C17").Value = "gross photosynthesis"
$ws.Range("D17").Value = 1600
$ws.Range("A18").Value = 4.9754327191513097
$ws.Range("A18").NumberFormat = "0"
$ws.Range("B18").Value = 22.3115577889447
$ws.Range("C18").Value = "gross photosynthesis"
$ws.Range("D18").Value = 1600
$ws.Range("A19").Value = 9.9871580122836399
$ws.Range("A19").NumberFormat = "0"
$ws.Range("B19").Value = 30.753768844221099
$ws.Range("C19").Value = "gross photosynthesis"
$ws.Range("D19").Value = 1600
$ws.Range("A20").Value = 14.9865996649916
$ws.Range("A20").NumberFormat = "0"
$ws.Range("B20").Value = 30.3517587939698
$ws.Range("C20").Value = "gross photosynthesis"
$ws.Range("D20").Value = 1600
$ws.Range("A21").Value = 20.039642657733101
$ws.Range("A21").NumberFormat = "0"
$ws.Range("B21").Value = 28.542713567839201
$ws.Range("C21").Value = "gross photosynthesis"
$ws.Range("D21").Value = 1600
$ws.Range("A22").Value = -0.099106644332778104
$ws.Range("A22").NumberFormat = "0"
$ws.Range("B22").Value = 8.6432160804019809
$ws.Range("C22").Value = "net photosynthesis"
$ws.Range("D22").Value = 1600
$ws.Range("A23").Value = 4.9614740368509196
$ws.Range("A23").NumberFormat = "0"
$ws.Range("B23").Value = 12.2613065326633
$ws.Range("C23").Value = "net photosynthesis"
$ws.Range("D23").Value = 1600
$ws.Range("A24").Value = 9.9648241206030104
$ws.Range("A24").NumberFormat = "0"
$ws.Range("B24").Value = 14.6733668341708
$ws.Range("C24").Value = "net photosynthesis"
$ws.Range("D24").Value = 1600
$ws.Range("A25").Value = 14.9972082635399
$ws.Range("A25").NumberFormat = "0"
$ws.Range("B25").Value = -2.0100502512562999
$ws.Range("C25").Value = "net photosynthesis"
$ws.Range("D25").Value = 1600
$ws.Range("A26").Value = 20.0530429927414
$ws.Range("A26").NumberFormat = "0"
$ws.Range("B26").Value = -1.80904522613067
$ws.Range("C26").Value = "net photosynthesis"
$ws.Range("D26").Value = 1600
$ws.Range("A27").Value = -0.010050251256281201
$ws.Range("A27").NumberFormat = "0"
$ws.Range("B27").Value = -7.2361809045226604
$ws.Range("C27").Value = "dark respiration"
$ws.Range("D27").Value = 1600
$ws.Range("A28").Value = 4.9302065884980397
$ws.Range("A28").NumberFormat = "0"
$ws.Range("B28").Value = -10.251256281407001
$ws.Range("C28").Value = "dark respiration"
$ws.Range("D28").Value = 1600
$ws.Range("A29").Value = 9.9221105527638098
$ws.Range("A29").NumberFormat = "0"
$ws.Range("B29").Value = -16.0804020100502
$ws.Range("C29").Value = "dark respiration"
$ws.Range("D29").Value = 1600
$ws.Range("A30").Value = 15.0106085985483
$ws.Range("A30").NumberFormat = "0"
$ws.Range("B30").Value = -32.361809045226202
$ws.Range("C30").Value = "dark respiration"
$ws.Range("D30").Value = 1600
$ws.Range("A31").Value = 20.013400335008299
$ws.Range("A31").NumberFormat = "0"
$ws.Range("B31").Value = -30.3517587939699
$ws.Range("C31").Value = "dark respiration"
$ws.Range("D31").Value = 1600

# --- Selection matches authored state ---
$ws.Range("B17:B31").Select()
